# Add a "Save" column (column H) to the s_vals sheet, mirroring the
# existing header style used by the other header cells (e.g. G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell
$ws.Range("H1").Value = "Save"

# Copy G1's formatting (bold, border, centered alignment) onto H1 so the
# new header matches the look of the existing headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
